# Atualizado por script em 24-11-2023 14:45
#
# Rwanda Premier League 2023-2024 sheet update:
#  - Row 57 and Row 58 swap match data (order of two matches played on the
#    same date/round was corrected).
#  - Rows 62-65 rotate match data by one position (62<-63, 63<-64, 64<-65,
#    65<-62), correcting the order of four matches played on the same date.
#  - A new match row (67) is appended for Gorilla vs Etoile de L'Est.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 57 and 58 (columns F:V) ---------------------------------
$ws.Range("F57").Value = 'Mukura Victory Sports'
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 'Etincelles'
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 1.7
$ws.Range("L57").Value = 1.65
$ws.Range("M57").Value = '28/10/2023 12:08'
$ws.Range("N57").Value = 3.15
$ws.Range("P57").Value = 3.27
$ws.Range("Q57").Value = '28/10/2023 13:02'
$ws.Range("R57").Value = 4.1
$ws.Range("T57").Value = 5.02
$ws.Range("U57").Value = '28/10/2023 12:08'
$ws.Range("V57").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/mukura-victory-sports-etincelles/hlvWrwMs/'

$ws.Range("F58").Value = 'Bugesera'
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 'Amagaju'
$ws.Range("I58").Value = 2
$ws.Range("J58").Value = 1.88
$ws.Range("L58").Value = 1.85
$ws.Range("M58").Value = '28/10/2023 14:10'
$ws.Range("N58").Value = 2.89
$ws.Range("P58").Value = 3.06
$ws.Range("Q58").Value = '28/10/2023 14:10'
$ws.Range("R58").Value = 3.69
$ws.Range("T58").Value = 4.16
$ws.Range("U58").Value = '28/10/2023 14:10'
$ws.Range("V58").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/bugesera-amagaju/juh4mJDQ/'

# --- Rotate rows 62-65 (columns F:V): 62<-63, 63<-64, 64<-65, 65<-62 ----
$ws.Range("F62").Value = 'Marines'
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 'Amagaju'
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2.16
$ws.Range("L62").Value = 2.02
$ws.Range("M62").Value = '04/11/2023 13:04'
$ws.Range("N62").Value = 2.86
$ws.Range("P62").Value = 2.91
$ws.Range("Q62").Value = '04/11/2023 13:04'
$ws.Range("R62").Value = 2.99
$ws.Range("T62").Value = 3.76
$ws.Range("U62").Value = '04/11/2023 13:04'
$ws.Range("V62").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/marines-amagaju/ANNcgotH/'

$ws.Range("F63").Value = 'Muhazi United'
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 'APR'
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = 5.61
$ws.Range("L63").Value = 4.23
$ws.Range("M63").Value = '04/11/2023 13:55'
$ws.Range("N63").Value = 3.67
$ws.Range("P63").Value = 3
$ws.Range("Q63").Value = '04/11/2023 13:55'
$ws.Range("R63").Value = 1.43
$ws.Range("T63").Value = 1.87
$ws.Range("U63").Value = '04/11/2023 13:55'
$ws.Range("V63").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/muhazi-united-apr/ryM5iPAT/'

$ws.Range("F64").Value = 'Musanze'
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 'Kiyovu'
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2.91
$ws.Range("L64").Value = 2.12
$ws.Range("M64").Value = '04/11/2023 13:13'
$ws.Range("N64").Value = 2.71
$ws.Range("P64").Value = 3.05
$ws.Range("Q64").Value = '04/11/2023 13:49'
$ws.Range("R64").Value = 2.3
$ws.Range("T64").Value = 3.09
$ws.Range("U64").Value = '04/11/2023 13:13'
$ws.Range("V64").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/musanze-kiyovu/OU3YnqJp/'

$ws.Range("F65").Value = 'Rayon Sport'
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 'Mukura Victory Sports'
$ws.Range("I65").Value = 1
$ws.Range("J65").Value = 1.71
$ws.Range("L65").Value = 1.72
$ws.Range("M65").Value = '04/11/2023 11:35'
$ws.Range("N65").Value = 3.03
$ws.Range("P65").Value = 3.15
$ws.Range("Q65").Value = '04/11/2023 12:02'
$ws.Range("R65").Value = 4.24
$ws.Range("T65").Value = 4.75
$ws.Range("U65").Value = '04/11/2023 11:35'
$ws.Range("V65").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/rayon-sport-mukura-victory-sports/4pN1h5eN/'

# --- Append new row 67 (Gorilla vs Etoile de L'Est) ---------------------
# Copy formatting from the last data row (66) so the new row's "Indice"
# (A) and "data_partida" (E) cells reuse the existing cell styles instead
# of creating new ones.
$ws.Range("A66").Copy() | Out-Null
$ws.Range("A67").PasteSpecial(-4122) | Out-Null
$ws.Range("E66").Copy() | Out-Null
$ws.Range("E67").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A67").Value = 66
$ws.Range("B67").Value = 'rwanda'
$ws.Range("C67").Value = 'premier-league'
$ws.Range("D67").Value = '2023-2024'
$ws.Range("E67").Value = 45254.58333333334
$ws.Range("F67").Value = 'Gorilla'
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = "Etoile de L'Est"
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2.19
$ws.Range("K67").Value = '10/11/2023 03:13'
$ws.Range("L67").Value = 2.17
$ws.Range("M67").Value = '24/11/2023 13:26'
$ws.Range("N67").Value = 2.85
$ws.Range("O67").Value = '10/11/2023 03:13'
$ws.Range("P67").Value = 2.76
$ws.Range("Q67").Value = '24/11/2023 13:26'
$ws.Range("R67").Value = 3.19
$ws.Range("S67").Value = '10/11/2023 03:13'
$ws.Range("T67").Value = 3.55
$ws.Range("U67").Value = '24/11/2023 13:26'
$ws.Range("V67").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/gorilla-etoile-de-l-est/pIXDgfz0/'
